# Weekly data refresh: insert a new week's record as row 322 and shift
# all subsequent rows (322-369) down by one (to 323-370).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 322, pushing the existing rows 322..369 down to 323..370.
$ws.Rows("322:322").Insert()

# Populate the newly inserted row 322 with this week's record.
$ws.Range("A322").Value = 3
$ws.Range("B322").Value = "Femacal de La Calera"
$ws.Range("C322").Value = "Coquimbo"
$ws.Range("D322").Value = 45173
$ws.Range("E322").Value = 5
$ws.Range("F322").Value = "Fruta"
$ws.Range("G322").Value = 100101
$ws.Range("H322").Value = "Berries"
$ws.Range("I322").Value = 100101001
$ws.Range("J322").Value = "Arándano (blue)"
$ws.Range("K322").Value = "Sin especificar"
$ws.Range("L322").Value = "Primera"
$ws.Range("M322").Value = 45
$ws.Range("N322").Value = 13000
$ws.Range("O322").Value = 13000
$ws.Range("P322").Value = 13000
$ws.Range("Q322").Value = "$/bandeja 12 canastillos 125 gramos"
$ws.Range("R322").Value = "Provincia de Quillota"
$ws.Range("S322").Value = 8667
$ws.Range("T322").Value = 1.5
